$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @("POST009112-5", "POST014122-5", "SET180632", "TCAT009203", "TCAT223010", "TCATOA223010", "TRA057", "TRA057出入", "TRA375")
$colB = @("秀林和平郵局(花蓮12支)", "和平梨山郵局(台中94支)", "統一超商  環山店", "統一速達    梨山衛星所", "統一速達  梨山衛星所", "統一速達  梨山衛星所", "台灣鐵路管理局　和平站", "台灣鐵路管理局　和平站  出入口", "台灣鐵路管理局  阿里山站")
$colC = @("花蓮縣秀林鄉和平村113號", "台中市和平區梨山村中正路89號", "台中市和平區中興路三段64-5、64-6號", "台中市和平區梨山村福壽路5號之1", "台中市和平區梨山村福壽路5號之1", "台中市和平區梨山村福壽路5號之1", "花蓮縣秀林鄉和平村２７６號", "花蓮縣秀林鄉和平村２７６號", "嘉義縣阿里山鄉中正村1號")

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}
